# Bump the "Förändrad" (Changed) date column (C) by one day (45189 -> 45190)
# for every data row in the sheet (rows 2 through 252).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 252
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)  # Column C
    if ($cell.Value2 -eq 45189) {
        $cell.Value2 = 45190
    }
}
